$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '49.856.12'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +3.67%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.647.99'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +5.96%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '113.80'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +7.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '327.11'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.21%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.30%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.12'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.18'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.55%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0821'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.22%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.34'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +4.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.064.78'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +6.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.643.17'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +5.87%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +4.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '49.806.80'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.17'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.61%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.93'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.30%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.18'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '276.68'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.81%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.76'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +3.85%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.03'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.19'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -4.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.02'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.96%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.30'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.58'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.23%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.55%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0808'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +4.05%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.07'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +7.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.76'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.85%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +6.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '125.10'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.85%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.56%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.21'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.04%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.84%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +4.64%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.074.39'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +3.53%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.27'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +13.18%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +5.83%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.13'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.69%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '81.69'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.07%  '
